# Generate Report for Handback
# Updates the localization-status workbook:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    on the Overview sheet (both language columns) and on each language sheet.
#  - Latest Handback DateTime is refreshed for zh-cn and de-de.
#  - The stale "handback file is not the latest" Error Detail is cleared now
#    that the handback is in sync.
#  - A couple of columns are widened / narrowed to fit the refreshed content.

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status columns for zh-cn (E2) and de-de (F2) ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-18 12:49:52"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-18 12:50:00"
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments ---
# Overview: widen the two "handoff/handback file" columns (E, F)
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# zh-cn: widen Status column (C), narrow Error Detail column (P)
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332

# de-de: widen Status column (C), narrow Error Detail column (P)
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
